$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the row above (row 11) into the new row 12 first
$ws.Range("A11:M11").Copy()
$ws.Range("A12:M12").PasteSpecial(-4122)

# Now set the values for row 12 (2021年)
$ws.Range("A12").Value = "2021年"
$ws.Range("B12").Value = 678330.1
$ws.Range("C12").Value = 6044545.3
$ws.Range("D12").Value = 20987578
$ws.Range("E12").Value = 51646571.5
$ws.Range("F12").Value = 35028188.4
$ws.Range("G12").Value = 20487092.6
$ws.Range("H12").Value = 60426739.7
$ws.Range("I12").Value = 2090164.2
$ws.Range("K12").Value = 645458.7
$ws.Range("L12").Value = 39439161.7
$ws.Range("M12").Value = 4233128.6
